# update version to v0.4 structure: target=common options, value=options by case
#
# Adds a column D holding the "common options" JSON blob (richly syntax
# colored, like a copy/paste from a JSON-highlighting editor) and rewrites
# C3 from a plain number into a similar rich-text JSON blob for the
# per-case "skip" option.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# colours used by the highlighted JSON text (VBA-style OLE RGB = R + G*256 + B*65536)
$purple = 9703559   # FF871094 - keys
$dark   = 526344    # FF080808 - punctuation / default text
$green  = 1539334   # FF067D17 - string values

# ---------------------------------------------------------------------
# D1: header cell, mirrors C1 ("click") both in value and style
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "click"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# D2: the "common options" JSON value -> {"target": "selector=#account .login"}
# ---------------------------------------------------------------------
$d2text = '{"target": "selector=#account .login"}'
$ws.Range("D2").Value = $d2text
$ws.Range("C1").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$ws.Range("D2").Characters(2,8).Font.Color = $purple    # "target"
$ws.Range("D2").Characters(10,2).Font.Color = $dark     # ": "
$ws.Range("D2").Characters(12,26).Font.Color = $green   # "selector=#account .login"
$ws.Range("D2").Characters(38,1).Font.Color = $dark     # }

# ---------------------------------------------------------------------
# D3: empty cell, only inherits C3's old plain-number style
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("D3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# C3: per-case option JSON value -> {"skip": true}
# ---------------------------------------------------------------------
$c3text = '{"skip": true}'
$ws.Range("C3").Value = $c3text
$ws.Range("C1").Copy()
$ws.Range("C3").PasteSpecial(-4122)

$ws.Range("C3").Characters(2,1).Font.Color = $purple                     # "
$ws.Range("C3").Characters(3,4).Font.Name = "ＭＳ Ｐゴシック"
$ws.Range("C3").Characters(3,4).Font.Color = $purple                     # skip
$ws.Range("C3").Characters(7,1).Font.Color = $purple                     # "
$ws.Range("C3").Characters(8,2).Font.Color = $dark                       # ": "
$ws.Range("C3").Characters(10,4).Font.Name = "ＭＳ Ｐゴシック"
$ws.Range("C3").Characters(10,4).Font.Color = $green                     # true
$ws.Range("C3").Characters(14,1).Font.Color = $dark                      # }

# ---------------------------------------------------------------------
# column D sizing (best effort - matches the author's AutoFit width)
# ---------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 41.93

# ---------------------------------------------------------------------
# selection, as left by the editing session
# ---------------------------------------------------------------------
$ws.Range("D6").Select()
